$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$win = $excel.ActiveWindow
Write-Host ("Initial ScrollColumn: " + $win.ScrollColumn())
Write-Host ("Initial ScrollRow: " + $win.ScrollRow())
$win.ScrollColumn = 20
$win.ScrollRow = 2
Write-Host ("New ScrollColumn: " + $win.ScrollColumn())
Write-Host ("New ScrollRow: " + $win.ScrollRow())
